$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 372.85715
$ws.Range("I18").Value = 351.66666
$ws.Range("K18").Value = 351.66666
$ws.Range("M18").Value = -67.66665999999998
$ws.Range("H41").Value = 38465416
$ws.Range("I41").Value = 1147.8
$ws.Range("K41").Value = 1147.8
$ws.Range("M41").Value = -707.8
$ws.Range("H96").Value = 834.6667
$ws.Range("I96").Value = 752.8333
$ws.Range("K96").Value = 2258.4999
$ws.Range("M96").Value = -885.4998999999998
$ws.Range("H98").Value = 3827.0715
$ws.Range("I98").Value = 3709.3333
$ws.Range("K98").Value = 3709.3333
$ws.Range("M98").Value = -2211.3333
$ws.Range("H107").Value = 32831.387
$ws.Range("I107").Value = 36291.93
$ws.Range("J107").Value = 533
$ws.Range("K107").Value = 36291.93
$ws.Range("L107").Value = 533
$ws.Range("M107").Value = -34371.93
$ws.Range("N107").Value = -4373
$ws.Range("H122").Value = 3827.0715
$ws.Range("I122").Value = 3709.3333
$ws.Range("K122").Value = 11127.9999
$ws.Range("M122").Value = -8677.999899999999
$ws.Range("H125").Value = 12349504
$ws.Range("J125").Value = 13892915
$ws.Range("L125").Value = 125036235
$ws.Range("N125").Value = -125041155
$ws.Range("H127").Value = 548
$ws.Range("I127").Value = 548
$ws.Range("K127").Value = 1644
$ws.Range("M127").Value = 3316
$ws.Range("H137").Value = 3800.7334
$ws.Range("I137").Value = 3690.4138
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 11071.2414
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -8521.241399999999
$ws.Range("N137").Value = -26100
$ws.Range("H138").Value = 4201.159
$ws.Range("I138").Value = 2808.111
$ws.Range("J138").Value = 4410.1167
$ws.Range("K138").Value = 8424.332999999999
$ws.Range("L138").Value = 13230.3501
$ws.Range("M138").Value = -3284.332999999999
$ws.Range("N138").Value = -23510.3501
$ws.Range("H141").Value = 3392.375
$ws.Range("I141").Value = 3014.4167
$ws.Range("K141").Value = 9043.250100000001
$ws.Range("M141").Value = -3863.250100000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 84988.5
$ws.Range("J123").Value = 84988.5
$ws.Range("L123").Value = 84988.5
$ws.Range("N123").Value = -94788.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 346484.78
$ws.Range("I107").Value = 1188.7826
$ws.Range("K107").Value = 1188.7826
$ws.Range("M107").Value = 731.2174
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40538.52
$ws.Range("I31").Value = 995.2727
$ws.Range("J31").Value = 67724.5
$ws.Range("K31").Value = 995.2727
$ws.Range("L31").Value = 67724.5
$ws.Range("M31").Value = -700.2727
$ws.Range("N31").Value = -68314.5
$ws.Range("H34").Value = 40538.52
$ws.Range("I34").Value = 995.2727
$ws.Range("J34").Value = 67724.5
$ws.Range("K34").Value = 995.2727
$ws.Range("L34").Value = 67724.5
$ws.Range("M34").Value = -793.2727
$ws.Range("N34").Value = -68128.5
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
$ws.Range("H99").Value = 4411
$ws.Range("I99").Value = 4362.2
$ws.Range("J99").Value = 4472
$ws.Range("K99").Value = 4362.2
$ws.Range("L99").Value = 4472
$ws.Range("M99").Value = -2864.2
$ws.Range("N99").Value = -7468
$ws.Range("H122").Value = 3088.85
$ws.Range("I122").Value = 2398.6667
$ws.Range("J122").Value = 3653.5454
$ws.Range("K122").Value = 7196.000100000001
$ws.Range("L122").Value = 10960.6362
$ws.Range("M122").Value = -4746.000100000001
$ws.Range("N122").Value = -15860.6362
$ws.Range("H126").Value = 4411
$ws.Range("I126").Value = 4362.2
$ws.Range("J126").Value = 4472
$ws.Range("K126").Value = 13086.6
$ws.Range("L126").Value = 13416
$ws.Range("M126").Value = -10616.6
$ws.Range("N126").Value = -18356
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 875
$ws.Range("I26").Value = 750
$ws.Range("K26").Value = 2250
$ws.Range("M26").Value = -1962
$ws.Range("H28").Value = 1000
$ws.Range("I28").Value = 1000
$ws.Range("K28").Value = 3000
$ws.Range("M28").Value = -2768
$ws.Range("H56").Value = 5579.6
$ws.Range("I56").Value = 5579.6
$ws.Range("K56").Value = 5579.6
$ws.Range("M56").Value = -5049.6
$ws.Range("H87").Value = 15532.5
$ws.Range("I87").Value = 15532.5
$ws.Range("K87").Value = 46597.5
$ws.Range("M87").Value = -45349.5
$ws.Range("H90").Value = 15532.5
$ws.Range("I90").Value = 15532.5
$ws.Range("K90").Value = 139792.5
$ws.Range("M90").Value = -133552.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 33679
$ws.Range("J57").Value = 48799.668
$ws.Range("L57").Value = 48799.668
$ws.Range("N57").Value = -50439.668
$ws.Range("H99").Value = 7192.75
$ws.Range("I99").Value = 7192.75
$ws.Range("K99").Value = 7192.75
$ws.Range("M99").Value = -4946.75
$ws.Range("H122").Value = 3150.524
$ws.Range("I122").Value = 2133.7144
$ws.Range("J122").Value = 5184.143
$ws.Range("K122").Value = 6401.1432
$ws.Range("L122").Value = 15552.429
$ws.Range("M122").Value = -3951.1432
$ws.Range("N122").Value = -20452.429
$ws.Range("H131").Value = 25162.5
$ws.Range("J131").Value = 25162.5
$ws.Range("L131").Value = 25162.5
$ws.Range("N131").Value = -35242.5
$ws.Range("H132").Value = 158929.58
$ws.Range("I132").Value = 19376.75
$ws.Range("K132").Value = 58130.25
$ws.Range("M132").Value = -55600.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8352.947
$ws.Range("I7").Value = 9308.154
$ws.Range("J7").Value = 6283.3335
$ws.Range("K7").Value = 9308.154
$ws.Range("L7").Value = 6283.3335
$ws.Range("M7").Value = -9196.154
$ws.Range("N7").Value = -6507.3335
$ws.Range("H22").Value = 611.0625
$ws.Range("I22").Value = 535.2727
$ws.Range("K22").Value = 535.2727
$ws.Range("M22").Value = -240.2727
$ws.Range("H27").Value = 611.0625
$ws.Range("I27").Value = 535.2727
$ws.Range("K27").Value = 535.2727
$ws.Range("M27").Value = -428.2727
$ws.Range("H68").Value = 112788.89
$ws.Range("I68").Value = 800
$ws.Range("K68").Value = 800
$ws.Range("M68").Value = -51
$ws.Range("H71").Value = 112788.89
$ws.Range("I71").Value = 800
$ws.Range("K71").Value = 4000
$ws.Range("M71").Value = -256
$ws.Range("H82").Value = 2849.8
$ws.Range("I82").Value = 1666
$ws.Range("J82").Value = 3357.1428
$ws.Range("K82").Value = 1666
$ws.Range("L82").Value = 3357.1428
$ws.Range("M82").Value = -1305
$ws.Range("N82").Value = -4079.1428
$ws.Range("H85").Value = 2849.8
$ws.Range("I85").Value = 1666
$ws.Range("J85").Value = 3357.1428
$ws.Range("K85").Value = 1666
$ws.Range("L85").Value = 3357.1428
$ws.Range("M85").Value = -418
$ws.Range("N85").Value = -5853.1428
$ws.Range("H122").Value = 3185.04
$ws.Range("I122").Value = 3029.7222
$ws.Range("J122").Value = 3584.4285
$ws.Range("K122").Value = 9089.1666
$ws.Range("L122").Value = 10753.2855
$ws.Range("M122").Value = -6639.1666
$ws.Range("N122").Value = -15653.2855
$ws.Range("H126").Value = 8352.947
$ws.Range("I126").Value = 9308.154
$ws.Range("J126").Value = 6283.3335
$ws.Range("K126").Value = 27924.462
$ws.Range("L126").Value = 18850.0005
$ws.Range("M126").Value = -25454.462
$ws.Range("N126").Value = -23790.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H74").Value = 8552.286
$ws.Range("I74").Value = 5999
$ws.Range("K74").Value = 5999
$ws.Range("M74").Value = -5063
$ws.Range("H77").Value = 8552.286
$ws.Range("I77").Value = 5999
$ws.Range("K77").Value = 17997
$ws.Range("M77").Value = -13317
$ws.Range("H107").Value = 761.64
$ws.Range("I107").Value = 801.5238000000001
$ws.Range("J107").Value = 552.25
$ws.Range("K107").Value = 2404.5714
$ws.Range("L107").Value = 1656.75
$ws.Range("M107").Value = -484.5714000000003
$ws.Range("N107").Value = -5496.75
$ws.Range("H122").Value = 27779168
$ws.Range("I122").Value = 33334804
$ws.Range("K122").Value = 100004412
$ws.Range("M122").Value = -100001962
$ws.Range("H126").Value = 722.6667
$ws.Range("I126").Value = 586.2857
$ws.Range("K126").Value = 1758.8571
$ws.Range("M126").Value = 711.1428999999998
